$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell data
$ws.Range('A1').Value = 'Name'
$ws.Range('B1').Value = 'Description'
$ws.Range('C1').Value = 'Scope'
$ws.Range('D1').Value = 'Limitations'
$ws.Range('E1').Value = 'Source'
$ws.Rows.Item(1).RowHeight = 17

$ws.Range('A2').Value = 'USDA Soybean Yield'
$ws.Range('B2').Value = 'Data on soybean yield and it''s various measures ie bushels total, bushels / acre, etc. '
$ws.Range('C2').Value = 'all US States, from 1996 to 2019'
$ws.Range('D2').Value = 'Missing vals'
$ws.Range('E2').Value = 'Quick stats lite query'
$ws.Rows.Item(2).RowHeight = 17

$ws.Range('A3').Value = 'Soybean Fertilizer'
$ws.Range('B3').Value = 'Fertilizer applied and stats'
$ws.Range('C3').Value = 'all US States, from 1996 to 2019'
$ws.Range('D3').Value = 'Missing vals'
$ws.Range('E3').Value = 'Quick stats lite query'
$ws.Rows.Item(3).RowHeight = 17

$ws.Range('A4').Value = 'Insecticides'
$ws.Range('B4').Value = 'Insecticides used and stats '
$ws.Range('C4').Value = 'all US States, from 1996 to 2019'
$ws.Range('D4').Value = 'Missing vals'
$ws.Range('E4').Value = 'Quick stats lite query'
$ws.Rows.Item(4).RowHeight = 17

$ws.Range('A5').Value = 'Fungicides'
$ws.Range('B5').Value = 'Fungicides used'
$ws.Range('C5').Value = 'all US States, from 1996 to 2019'
$ws.Range('D5').Value = 'Missing vals'
$ws.Range('E5').Value = 'Quick stats lite query'
$ws.Rows.Item(5).RowHeight = 17

$ws.Range('A6').Value = 'Herbicides'
$ws.Range('B6').Value = 'Herbicides used that year'
$ws.Range('C6').Value = 'all US States, from 1996 to 2019'
$ws.Range('D6').Value = 'Missing vals'
$ws.Range('E6').Value = 'Quick stats lite query'
$ws.Rows.Item(6).RowHeight = 17

$ws.Range('A7').Value = 'Average temperature'
$ws.Range('B7').Value = 'Average of the daily average termperature observed over a month'
$ws.Range('C7').Value = 'all US States, from 1996 to 2019'
$ws.Range('E7').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(7).RowHeight = 17

$ws.Range('A8').Value = 'Minimum temperature'
$ws.Range('B8').Value = 'Average daily minimum temperature observed over a month'
$ws.Range('C8').Value = 'all US States, from 1996 to 2019'
$ws.Range('E8').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(8).RowHeight = 17

$ws.Range('A9').Value = 'Maximum temperature'
$ws.Range('B9').Value = 'Average daily maximum temperature observed over a month'
$ws.Range('C9').Value = 'all US States, from 1996 to 2019'
$ws.Range('E9').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(9).RowHeight = 17

$ws.Range('A10').Value = 'Cooling degree days'
$ws.Range('B10').Value = 'It is the measurement of number of degrees that a day''s average temperature is above 65 degree F (annual statistic)'
$ws.Range('C10').Value = 'all US States, from 1996 to 2019'
$ws.Range('E10').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(10).RowHeight = 34

$ws.Range('A11').Value = 'Heating degree days'
$ws.Range('B11').Value = 'It is the measurement of number of degrees that a day''s average temperature is below 65 degree F (annual statistic)'
$ws.Range('C11').Value = 'all US States, from 1996 to 2019'
$ws.Range('E11').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(11).RowHeight = 34

$ws.Range('A12').Value = 'Palmer Drought Severity Index (PDSI)'
$ws.Range('B12').Value = 'Measures the duration and internsity of droughts by measuring the dryness based on precipitation and temperature data as well as the local Available Water Content of the soil (Monthly statistic)'
$ws.Range('C12').Value = 'all US States, from 1996 to 2019'
$ws.Range('E12').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(12).RowHeight = 51

$ws.Range('A13').Value = 'Palmer Hydrological Drought Index (PHDI)'
$ws.Range('B13').Value = 'Measures hydrological impacts of drought (e.g., reservoir levels, groundwater levels, etc.) which take longer to develop and longer to recover from.(Monthly statistic)'
$ws.Range('C13').Value = 'all US States, from 1996 to 2019'
$ws.Range('E13').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(13).RowHeight = 34

$ws.Range('A14').Value = 'Palmer z-index'
$ws.Range('B14').Value = 'Measures the short-term drought on a monthly scale'
$ws.Range('C14').Value = 'all US States, from 1996 to 2019'
$ws.Range('E14').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(14).RowHeight = 17

$ws.Range('A15').Value = 'Palmer Modified Drought Index (PMDI)'
$ws.Range('B15').Value = 'Operational version of PDSI'
$ws.Range('C15').Value = 'all US States, from 1996 to 2019'
$ws.Range('E15').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(15).RowHeight = 17

$ws.Range('A16').Value = 'Soybeans Crop moisture stress index (CMSI)'
$ws.Range('B16').Value = 'Measure the impact of both lack and abundance of soil moisture on the National crop yield of Soybeans. It is calculated using the Palmer Z index and Annual average crop productitivy values within each US climate division'
$ws.Range('C16').Value = 'all US States, from 1996 to 2019'
$ws.Range('E16').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(16).RowHeight = 51

$ws.Range('A17').Value = 'Precipitation '
$ws.Range('B17').Value = 'Average precipitation over a month'
$ws.Range('C17').Value = 'all US States, from 1996 to 2019'
$ws.Range('E17').Value = 'Climate at a glance (NOAA)'
$ws.Rows.Item(17).RowHeight = 17

# Column widths and wrap text (applied after data so new rows inherit formatting)
$ws.Columns.Item(1).ColumnWidth = 35.6640625
$ws.Columns.Item(2).ColumnWidth = 72.1640625
$ws.Range("B1:B17").WrapText = $true

# Selection
$ws.Range('B17').Select() | Out-Null
